$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Jurisdiction codes (column A) and updated unique_customer_count values (column B)
# Written in ascending alphabetical order of the ISO jurisdiction code, row by row,
# starting at row 2 (row 1 holds the headers "jurisdiction" / "unique_customer_count").
$ws.Range("A2").Value = "AD"
$ws.Range("B2").Value = 4
$ws.Range("A3").Value = "AE"
$ws.Range("B3").Value = 4
$ws.Range("A4").Value = "AF"
$ws.Range("B4").Value = 21
$ws.Range("A5").Value = "AL"
$ws.Range("B5").Value = 8
$ws.Range("A6").Value = "AO"
$ws.Range("B6").Value = 0
$ws.Range("A7").Value = "AQ"
$ws.Range("B7").Value = 1
$ws.Range("A8").Value = "AS"
$ws.Range("B8").Value = 1
$ws.Range("A9").Value = "AT"
$ws.Range("B9").Value = 0
$ws.Range("A10").Value = "AU"
$ws.Range("B10").Value = 12
$ws.Range("A11").Value = "AX"
$ws.Range("B11").Value = 7
$ws.Range("A12").Value = "BE"
$ws.Range("B12").Value = 4
$ws.Range("A13").Value = "BF"
$ws.Range("B13").Value = 1
$ws.Range("A14").Value = "BM"
$ws.Range("B14").Value = 0
$ws.Range("A15").Value = "BT"
$ws.Range("B15").Value = 1
$ws.Range("A16").Value = "CA"
$ws.Range("B16").Value = 0
$ws.Range("A17").Value = "CH"
$ws.Range("B17").Value = 1
$ws.Range("A18").Value = "CL"
$ws.Range("B18").Value = 1
$ws.Range("A19").Value = "CM"
$ws.Range("B19").Value = 1
$ws.Range("A20").Value = "CN"
$ws.Range("B20").Value = 29
$ws.Range("A21").Value = "CO"
$ws.Range("B21").Value = 1
$ws.Range("A22").Value = "CU"
$ws.Range("B22").Value = 1
$ws.Range("A23").Value = "CW"
$ws.Range("B23").Value = 1
$ws.Range("A24").Value = "CY"
$ws.Range("B24").Value = 0
$ws.Range("A25").Value = "DE"
$ws.Range("B25").Value = 3
$ws.Range("A26").Value = "DM"
$ws.Range("B26").Value = 1
$ws.Range("A27").Value = "DO"
$ws.Range("B27").Value = 0
$ws.Range("A28").Value = "DZ"
$ws.Range("B28").Value = 14
$ws.Range("A29").Value = "EE"
$ws.Range("B29").Value = 0
$ws.Range("A30").Value = "ES"
$ws.Range("B30").Value = 1
$ws.Range("A31").Value = "FR"
$ws.Range("B31").Value = 5
$ws.Range("A32").Value = "GB"
$ws.Range("B32").Value = 65
$ws.Range("A33").Value = "GE"
$ws.Range("B33").Value = 0
$ws.Range("A34").Value = "GF"
$ws.Range("B34").Value = 0
$ws.Range("A35").Value = "GN"
$ws.Range("B35").Value = 1
$ws.Range("A36").Value = "GR"
$ws.Range("B36").Value = 0
$ws.Range("A37").Value = "GT"
$ws.Range("B37").Value = 1
$ws.Range("A38").Value = "GW"
$ws.Range("B38").Value = 1
$ws.Range("A39").Value = "HK"
$ws.Range("B39").Value = 1
$ws.Range("A40").Value = "HN"
$ws.Range("B40").Value = 0
$ws.Range("A41").Value = "HU"
$ws.Range("B41").Value = 0
$ws.Range("A42").Value = "ID"
$ws.Range("B42").Value = 1
$ws.Range("A43").Value = "IL"
$ws.Range("B43").Value = 0
$ws.Range("A44").Value = "IN"
$ws.Range("B44").Value = 25
$ws.Range("A45").Value = "IQ"
$ws.Range("B45").Value = 0
$ws.Range("A46").Value = "IR"
$ws.Range("B46").Value = 2
$ws.Range("A47").Value = "IT"
$ws.Range("B47").Value = 0
$ws.Range("A48").Value = "JE"
$ws.Range("B48").Value = 0
$ws.Range("A49").Value = "JP"
$ws.Range("B49").Value = 2
$ws.Range("A50").Value = "KP"
$ws.Range("B50").Value = 0
$ws.Range("A51").Value = "KR"
$ws.Range("B51").Value = 1
$ws.Range("A52").Value = "KW"
$ws.Range("B52").Value = 0
$ws.Range("A53").Value = "KY"
$ws.Range("B53").Value = 1
$ws.Range("A54").Value = "LS"
$ws.Range("B54").Value = 0
$ws.Range("A55").Value = "LU"
$ws.Range("B55").Value = 4
$ws.Range("A56").Value = "LV"
$ws.Range("B56").Value = 0
$ws.Range("A57").Value = "LY"
$ws.Range("B57").Value = 1
$ws.Range("A58").Value = "MO"
$ws.Range("B58").Value = 0
$ws.Range("A59").Value = "MV"
$ws.Range("B59").Value = 1
$ws.Range("A60").Value = "MW"
$ws.Range("B60").Value = 0
$ws.Range("A61").Value = "NF"
$ws.Range("B61").Value = 1
$ws.Range("A62").Value = "NL"
$ws.Range("B62").Value = 1
$ws.Range("A63").Value = "NO"
$ws.Range("B63").Value = 1
$ws.Range("A64").Value = "RU"
$ws.Range("B64").Value = 2
$ws.Range("A65").Value = "SE"
$ws.Range("B65").Value = 0
$ws.Range("A66").Value = "SG"
$ws.Range("B66").Value = 1
$ws.Range("A67").Value = "TH"
$ws.Range("B67").Value = 0
$ws.Range("A68").Value = "TN"
$ws.Range("B68").Value = 0
$ws.Range("A69").Value = "TO"
$ws.Range("B69").Value = 0
$ws.Range("A70").Value = "UG"
$ws.Range("B70").Value = 2
$ws.Range("A71").Value = "US"
$ws.Range("B71").Value = 16
$ws.Range("A72").Value = "UZ"
$ws.Range("B72").Value = 1
$ws.Range("A73").Value = "VA"
$ws.Range("B73").Value = 0
$ws.Range("A74").Value = "VN"
$ws.Range("B74").Value = 0
$ws.Range("A75").Value = "YE"
$ws.Range("B75").Value = 1
$ws.Range("A76").Value = "ZM"
$ws.Range("B76").Value = 0
$ws.Range("A77").Value = "ZW"
$ws.Range("B77").Value = 0
